$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.938.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.095.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.74"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.84%  "

$ws.Range("E10").Value = "  +2.36%  "

$ws.Range("E11").Value = "  +2.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.392.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.46"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.766"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.095.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.883.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("E19").Value = "  -3.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "

$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.93%  "

$ws.Range("E27").Value = "  +11.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.07%  "

$ws.Range("E31").Value = "  +0.65%  "

$ws.Range("E32").Value = "  +4.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.11%  "

$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("E35").Value = "  +1.51%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.30%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.62%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.45"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.100"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.27%  "

$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.455.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("E46").Value = "  +3.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.19%  "

$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.45%  "

$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.61%  "

$ws.Range("E50").Value = "  +2.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.287.63"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.79%  "
